$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.119.68"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").Value = "3.085.32"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'553.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").Value = "'137.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "3.079.75"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "'0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").Value = "'6.66"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.98%  "

$ws.Range("D12").Value = "'0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").Value = "'35.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "'0.0000218"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "3.581.22"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Value = "63.182.50"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "3.096.79"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "'502.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "

$ws.Range("D20").Value = "'6.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.65%  "

$ws.Range("D21").Value = "'13.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.79%  "

$ws.Range("D23").Value = "'7.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.38%  "

$ws.Range("D24").Value = "'78.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").Value = "'12.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Value = "'2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").Value = "'8.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("D29").Value = "'2.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "'26.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.35%  "

$ws.Range("D32").Value = "'2.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").Value = "'1.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.55%  "

$ws.Range("D34").Value = "'59.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.72%  "

$ws.Range("D35").Value = "'530.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.15%  "

$ws.Range("D36").Value = "'5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "'5.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").Value = "'0.0412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.59%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.073.81"
$ws.Range("E39").Value = "  +2.66%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0794"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("D42").Value = "'8.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

$ws.Range("D43").Value = "'2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.05%  "

$ws.Range("D44").Value = "'0.254"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.71%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "'2.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'120.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "

$ws.Range("D48").Value = "'23.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.95%  "

$ws.Range("D49").Value = "'0.106"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "'2.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +62.07%  "

$ws.Range("B51").Value = "PEPE"
$ws.Range("C51").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D51").Value = "0.0₃0496"
$ws.Range("E51").Value = "  -4.93%  "
